# Commit #4 "add some data": append 4 new missile rows (arrowred, arrowlight,
# bluepea, greenpea) to the "Missile" sheet / "表1" table, fix up the
# FrameTime value for the existing "dragonball" row, grow the table/
# worksheet dimensions accordingly, and leave the selection where the
# author's cursor ended up (E11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix existing data: dragonball's FrameTime was corrected from 3 to 6 ---
$ws.Range("G8").Value = 6

# --- new row 9: arrowred / 火箭 ---
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "arrowred"
$ws.Range("C9").Value = "火箭"
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = 50
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1

# --- new row 10: arrowlight / 光箭 ---
$ws.Range("A10").Value = 6
$ws.Range("B10").Value = "arrowlight"
$ws.Range("C10").Value = "光箭"
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 60
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1

# --- new row 11: bluepea / 蓝色豆子 ---
$ws.Range("A11").Value = 7
$ws.Range("B11").Value = "bluepea"
$ws.Range("C11").Value = "蓝色豆子"
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 70
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1

# --- new row 12: greenpea / 绿色豆子 ---
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = "greenpea"
$ws.Range("C12").Value = "绿色豆子"
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 80
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1

# --- grow the table ("表1") / autofilter so it covers the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:G12"))

# --- leave the active selection where the author's last edit landed ---
$ws.Range("E11").Select()
